# Update "想去人数" (interested/want-to-go counts) figures in column F
# across the four sheets, reflecting the data refresh captured in the
# commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 37
$ws.Range("F7").Value = 760
$ws.Range("F12").Value = 222
$ws.Range("F13").Value = 93
$ws.Range("F14").Value = 907
$ws.Range("F16").Value = 2025
$ws.Range("F17").Value = 516
$ws.Range("F18").Value = 8229
$ws.Range("F19").Value = 692

# 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 4
$ws.Range("F10").Value = 2

# 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5572
$ws.Range("F3").Value = 412
$ws.Range("F4").Value = 400

# 全部类型 (All types - aggregate of the above)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5572
$ws.Range("F4").Value = 412
$ws.Range("F5").Value = 400
$ws.Range("F8").Value = 4
$ws.Range("F11").Value = 37
$ws.Range("F13").Value = 760
$ws.Range("F20").Value = 222
$ws.Range("F22").Value = 93
$ws.Range("F24").Value = 907
$ws.Range("F27").Value = 2
$ws.Range("F28").Value = 2026
$ws.Range("F29").Value = 516
$ws.Range("F30").Value = 8229
$ws.Range("F33").Value = 692
